$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update activation date (keep as text, not a date serial)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Add English translation of Objectives (row 11)
$objText = "To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues."
$ws.Range("B11").Value = $objText
$ws.Range("C11").Value = $objText
$ws.Range("B11").Style = $ws.Range("B10").Style
$ws.Range("C11").Style = $ws.Range("C10").Style

# Add English translation of Short syllabus (row 15)
$shortSyllabus = "The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes"
$ws.Range("B15").Value = $shortSyllabus
$ws.Range("C15").Value = $shortSyllabus
$ws.Range("B15").Style = $ws.Range("B14").Style
$ws.Range("C15").Style = $ws.Range("C14").Style

# Add English translation of Syllabus (row 17)
$syllabus = "- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment."
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus
$ws.Range("B17").Style = $ws.Range("B16").Style
$ws.Range("C17").Style = $ws.Range("C16").Style

# Replace bibliography text
$bib = "LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014."
$ws.Range("B22").Value = $bib
$ws.Range("C22").Value = $bib
